# Auto-generated Excel COM-interop script to update Mateus_Profits market-data sheets
# (values refreshed by the scheduled market-data runner; no formulas involved).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 13162534
$ws.Range("I40").Value = 3300
$ws.Range("K40").Value = 3300
$ws.Range("M40").Value = -3125
$ws.Range("H51").Value = 7165
$ws.Range("J51").Value = 3995
$ws.Range("L51").Value = 3995
$ws.Range("N51").Value = -4963
$ws.Range("H58").Value = 1430.1428
$ws.Range("J58").Value = 1619.3334
$ws.Range("L58").Value = 4858.0002
$ws.Range("N58").Value = -5158.0002
$ws.Range("H80").Value = 12361.471
$ws.Range("I80").Value = 9477.637000000001
$ws.Range("K80").Value = 28432.911
$ws.Range("M80").Value = -27434.911
$ws.Range("H83").Value = 12361.471
$ws.Range("I83").Value = 9477.637000000001
$ws.Range("K83").Value = 85298.73300000001
$ws.Range("M83").Value = -80306.73300000001
$ws.Range("H121").Value = 1264.56
$ws.Range("J121").Value = 1274.6522
$ws.Range("L121").Value = 3823.9566
$ws.Range("N121").Value = -7317.9566
$ws.Range("H137").Value = 9640.429
$ws.Range("I137").Value = 12662.375
$ws.Range("K137").Value = 37987.125
$ws.Range("M137").Value = -35437.125
$ws.Range("H138").Value = 5339.0586
$ws.Range("I138").Value = 6808.6665
$ws.Range("K138").Value = 20425.9995
$ws.Range("M138").Value = -15285.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 49.5
$ws.Range("I5").Value = 49.5
$ws.Range("K5").Value = 49.5
$ws.Range("M5").Value = 62.5
$ws.Range("H45").Value = 3842.0286
$ws.Range("I45").Value = 2981.1904
$ws.Range("J45").Value = 5133.2856
$ws.Range("K45").Value = 2981.1904
$ws.Range("L45").Value = 5133.2856
$ws.Range("M45").Value = -2604.1904
$ws.Range("N45").Value = -5887.2856
$ws.Range("H61").Value = 11064.0625
$ws.Range("I61").Value = 11064.0625
$ws.Range("K61").Value = 11064.0625
$ws.Range("M61").Value = -10852.0625
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H74").Value = 4380.75
$ws.Range("I74").Value = 2787.125
$ws.Range("K74").Value = 2787.125
$ws.Range("M74").Value = -1913.125
$ws.Range("H77").Value = 4380.75
$ws.Range("I77").Value = 2787.125
$ws.Range("K77").Value = 13935.625
$ws.Range("M77").Value = -9567.625
$ws.Range("H133").Value = 79996
$ws.Range("J133").Value = 79996
$ws.Range("L133").Value = 79996
$ws.Range("N133").Value = -85056
$ws.Range("H134").Value = 108806.664
$ws.Range("J134").Value = 108806.664
$ws.Range("L134").Value = 108806.664
$ws.Range("N134").Value = -118946.664
$ws.Range("H136").Value = 11064.0625
$ws.Range("I136").Value = 11064.0625
$ws.Range("K136").Value = 33192.1875
$ws.Range("M136").Value = -30642.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 49.5
$ws.Range("I4").Value = 49.5
$ws.Range("K4").Value = 49.5
$ws.Range("M4").Value = 65.5
$ws.Range("H105").Value = 1406.6
$ws.Range("I105").Value = 1406.6
$ws.Range("K105").Value = 1406.6
$ws.Range("M105").Value = 340.4000000000001
$ws.Range("H134").Value = 1374.1818
$ws.Range("I134").Value = 1374.1818
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4122.5454
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1587.5454
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6022.579
$ws.Range("I31").Value = 4043.75
$ws.Range("J31").Value = 9414.857
$ws.Range("K31").Value = 4043.75
$ws.Range("L31").Value = 9414.857
$ws.Range("M31").Value = -3748.75
$ws.Range("N31").Value = -10004.857
$ws.Range("H34").Value = 6022.579
$ws.Range("I34").Value = 4043.75
$ws.Range("J34").Value = 9414.857
$ws.Range("K34").Value = 4043.75
$ws.Range("L34").Value = 9414.857
$ws.Range("M34").Value = -3841.75
$ws.Range("N34").Value = -9818.857
$ws.Range("H58").Value = 5249.2085
$ws.Range("I58").Value = 2034.1765
$ws.Range("K58").Value = 2034.1765
$ws.Range("M58").Value = -1831.1765
$ws.Range("H136").Value = 5249.2085
$ws.Range("I136").Value = 2034.1765
$ws.Range("K136").Value = 6102.529500000001
$ws.Range("M136").Value = -3552.529500000001
$ws.Range("H140").Value = 393083.88
$ws.Range("J140").Value = 393083.88
$ws.Range("L140").Value = 393083.88
$ws.Range("N140").Value = -403443.88
$ws.Range("H141").Value = 214839.8
$ws.Range("J141").Value = 259424.33
$ws.Range("L141").Value = 259424.33
$ws.Range("N141").Value = -269784.33

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 31340106
$ws.Range("I4").Value = 30645422
$ws.Range("K4").Value = 91936266
$ws.Range("M4").Value = -91936154
$ws.Range("H34").Value = 900
$ws.Range("H103").Value = 343.1
$ws.Range("I103").Value = 372.8
$ws.Range("J103").Value = 313.4
$ws.Range("K103").Value = 1118.4
$ws.Range("L103").Value = 940.1999999999999
$ws.Range("M103").Value = -239.4000000000001
$ws.Range("N103").Value = -2698.2
$ws.Range("H109").Value = 2198.9092
$ws.Range("I109").Value = 173.5
$ws.Range("J109").Value = 7600
$ws.Range("K109").Value = 520.5
$ws.Range("L109").Value = 22800
$ws.Range("M109").Value = 519.5
$ws.Range("N109").Value = -24880
$ws.Range("H112").Value = 10326
$ws.Range("I112").Value = 10326
$ws.Range("K112").Value = 30978
$ws.Range("M112").Value = -29870
$ws.Range("H120").Value = 17399.8
$ws.Range("I120").Value = 13499.5
$ws.Range("K120").Value = 40498.5
$ws.Range("M120").Value = -35660.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5399.4287
$ws.Range("I126").Value = 4553.8184
$ws.Range("K126").Value = 13661.4552
$ws.Range("M126").Value = -11191.4552

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 3612.25
$ws.Range("J12").Value = 3612.25
$ws.Range("L12").Value = 3612.25
$ws.Range("N12").Value = -3952.25
$ws.Range("H16").Value = 411.5
$ws.Range("J16").Value = 691.75
$ws.Range("L16").Value = 691.75
$ws.Range("N16").Value = -1031.75
$ws.Range("H22").Value = 3631.875
$ws.Range("I22").Value = 2714.75
$ws.Range("J22").Value = 4549
$ws.Range("K22").Value = 2714.75
$ws.Range("L22").Value = 4549
$ws.Range("M22").Value = -2419.75
$ws.Range("N22").Value = -5139
$ws.Range("H27").Value = 3631.875
$ws.Range("I27").Value = 2714.75
$ws.Range("J27").Value = 4549
$ws.Range("K27").Value = 2714.75
$ws.Range("L27").Value = 4549
$ws.Range("M27").Value = -2607.75
$ws.Range("N27").Value = -4763
$ws.Range("H35").Value = 1120
$ws.Range("I35").Value = 1120
$ws.Range("K35").Value = 1120
$ws.Range("M35").Value = -784
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H58").Value = 29000
$ws.Range("J58").Value = 29000
$ws.Range("L58").Value = 29000
$ws.Range("N58").Value = -29520

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 78735.5
$ws.Range("J46").Value = 80185.38
$ws.Range("L46").Value = 80185.38
$ws.Range("N46").Value = -80647.38
$ws.Range("H81").Value = 3991.0557
$ws.Range("I81").Value = 3353.6875
$ws.Range("K81").Value = 6707.375
$ws.Range("M81").Value = -5646.375
$ws.Range("H84").Value = 3991.0557
$ws.Range("I84").Value = 3353.6875
$ws.Range("K84").Value = 33536.875
$ws.Range("M84").Value = -28232.875
$ws.Range("H113").Value = 2185.1904
$ws.Range("J113").Value = 5103.8
$ws.Range("L113").Value = 15311.4
$ws.Range("N113").Value = -19651.4
$ws.Range("H134").Value = 78735.5
$ws.Range("J134").Value = 80185.38
$ws.Range("L134").Value = 240556.14
$ws.Range("N134").Value = -245626.14
